# Applies the cryptos-list refresh described by the commit:
#   "Updated cryptos list on Tue Oct 10 03:20:17 UTC 2023 with GitHub Actions"
#
# Every target cell holds plain text (prices/links/percentages are stored as
# strings, not numbers/formulas), so each update is a simple literal replace.
# The new value is written with a leading apostrophe so Excel treats it as
# literal text instead of silently coercing numeric-looking strings (e.g.
# "207.50" or "  -0.38%  ") into numbers/percentages, which would both change
# the stored type and mangle formatting (dropped zeros, trimmed spaces). The
# apostrophe-prefix entry tags the cell with Excel's built-in quote-prefix
# style, so Style is reset to 'Normal' right after to keep the cell's style
# index exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '27.709.13' },
    @{ Cell = 'E2'; Value = '  -0.38%  ' },
    @{ Cell = 'D3'; Value = '1.590.25' },
    @{ Cell = 'E3'; Value = '  -2.23%  ' },
    @{ Cell = 'E4'; Value = '  +0.65%  ' },
    @{ Cell = 'D5'; Value = '207.47' },
    @{ Cell = 'E5'; Value = '  -1.64%  ' },
    @{ Cell = 'E6'; Value = '  -2.83%  ' },
    @{ Cell = 'E7'; Value = '  +0.69%  ' },
    @{ Cell = 'D8'; Value = '22.25' },
    @{ Cell = 'E8'; Value = '  -3.99%  ' },
    @{ Cell = 'D9'; Value = '0.252' },
    @{ Cell = 'E9'; Value = '  -1.46%  ' },
    @{ Cell = 'E10'; Value = '  -2.52%  ' },
    @{ Cell = 'D11'; Value = '0.0869' },
    @{ Cell = 'E11'; Value = '  -1.02%  ' },
    @{ Cell = 'D12'; Value = '1.815.62' },
    @{ Cell = 'E12'; Value = '  -2.24%  ' },
    @{ Cell = 'D13'; Value = '1.570.39' },
    @{ Cell = 'E13'; Value = '  -3.49%  ' },
    @{ Cell = 'D14'; Value = '3.87' },
    @{ Cell = 'E14'; Value = '  -3.53%  ' },
    @{ Cell = 'E15'; Value = '  -4.50%  ' },
    @{ Cell = 'D16'; Value = '63.58' },
    @{ Cell = 'E16'; Value = '  -1.95%  ' },
    @{ Cell = 'D17'; Value = '27.690.43' },
    @{ Cell = 'E17'; Value = '  -0.51%  ' },
    @{ Cell = 'D18'; Value = '220.42' },
    @{ Cell = 'E18'; Value = '  -3.39%  ' },
    @{ Cell = 'E19'; Value = '  -2.92%  ' },
    @{ Cell = 'D20'; Value = '7.35' },
    @{ Cell = 'E20'; Value = '  -3.53%  ' },
    @{ Cell = 'E21'; Value = '  +0.65%  ' },
    @{ Cell = 'E22'; Value = '  -4.48%  ' },
    @{ Cell = 'D23'; Value = '9.61' },
    @{ Cell = 'E23'; Value = '  -3.09%  ' },
    @{ Cell = 'E24'; Value = '  -3.75%  ' },
    @{ Cell = 'D25'; Value = '153.35' },
    @{ Cell = 'E25'; Value = '  -1.16%  ' },
    @{ Cell = 'D26'; Value = '6.88' },
    @{ Cell = 'E26'; Value = '  -0.75%  ' },
    @{ Cell = 'E27'; Value = '  +0.69%  ' },
    @{ Cell = 'D28'; Value = '15.15' },
    @{ Cell = 'E28'; Value = '  -1.83%  ' },
    @{ Cell = 'E29'; Value = '  -4.14%  ' },
    @{ Cell = 'D30'; Value = '1.15' },
    @{ Cell = 'E30'; Value = '  -1.89%  ' },
    @{ Cell = 'D31'; Value = '0.0470' },
    @{ Cell = 'E31'; Value = '  -2.06%  ' },
    @{ Cell = 'E32'; Value = '  -5.07%  ' },
    @{ Cell = 'D33'; Value = '1.371.78' },
    @{ Cell = 'E33'; Value = '  -2.47%  ' },
    @{ Cell = 'D34'; Value = '2.93' },
    @{ Cell = 'E34'; Value = '  -5.24%  ' },
    @{ Cell = 'D35'; Value = '1.54' },
    @{ Cell = 'E35'; Value = '  -4.11%  ' },
    @{ Cell = 'D36'; Value = '0.980' },
    @{ Cell = 'E36'; Value = '  -1.75%  ' },
    @{ Cell = 'E37'; Value = '  -0.64%  ' },
    @{ Cell = 'E38'; Value = '  -0.90%  ' },
    @{ Cell = 'D39'; Value = '0.540' },
    @{ Cell = 'E39'; Value = '  -2.33%  ' },
    @{ Cell = 'D40'; Value = '0.826' },
    @{ Cell = 'E40'; Value = '  -2.56%  ' },
    @{ Cell = 'E41'; Value = '  +0.64%  ' },
    @{ Cell = 'D42'; Value = '0.969' },
    @{ Cell = 'E42'; Value = '  -3.13%  ' },
    @{ Cell = 'D43'; Value = '64.34' },
    @{ Cell = 'E43'; Value = '  -2.15%  ' },
    @{ Cell = 'D44'; Value = '2.17' },
    @{ Cell = 'E44'; Value = '  +2.81%  ' },
    @{ Cell = 'D45'; Value = '5.25' },
    @{ Cell = 'E45'; Value = '  -2.98%  ' },
    @{ Cell = 'E46'; Value = '  -3.93%  ' },
    @{ Cell = 'D47'; Value = '1.725.63' },
    @{ Cell = 'E47'; Value = '  -2.28%  ' },
    @{ Cell = 'D48'; Value = '87.83' },
    @{ Cell = 'E48'; Value = '  -0.58%  ' },
    @{ Cell = 'B49'; Value = 'BabyDogeCoin' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' },
    @{ Cell = 'D49'; Value = '0.0₆0100' },
    @{ Cell = 'E49'; Value = '  +13.01%  ' },
    @{ Cell = 'B50'; Value = 'Algorand' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'D50'; Value = '0.0970' },
    @{ Cell = 'E50'; Value = '  -3.90%  ' },
    @{ Cell = 'E51'; Value = '  -1.04%  ' }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.Value = "'" + $u.Value
    $r.Style = 'Normal'
}

